# NYPD 101st Precinct CompStat weekly report refresh:
#  - bump the "Volume NN  Number NN" header to the new issue number
#  - roll the "Report Covering the Week ... Through ..." date range forward one week
#  - refresh the crime-stat grid (rows 15-30, cols C:N) with the newly collected figures

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header text tweaks -------------------------------------------------
# A8 = "Volume 31   Number  51"  ->  "...  52"
$ws.Range("A8").Characters(21, 2).Text = "52"

# C9 = "Report Covering the Week  12/16/2024  Through  12/22/2024"
$ws.Range("C9").Characters(27, 10).Text = "12/23/2024"
$ws.Range("C9").Characters(48, 10).Text = "12/29/2024"

# ---- Helper: stable donor cells for the two "no data" placeholder styles
#      (C30 = style 13 / shared-string "0", E30 = style 13 / shared-string "***.*")
$zeroDonor = $ws.Range("C30")
$naDonor   = $ws.Range("E30")

function Set-Placeholder($cell, $donor) {
    $donor.Copy($cell)
}

# ---- Row 15 (Rape) -------------------------------------------------------
Set-Placeholder $ws.Range("C15") $zeroDonor
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 200

# ---- Row 16 (Robbery) -----------------------------------------------------
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 0
$ws.Range("G16").Value = 4
$ws.Range("H16").Value = 125
$ws.Range("I16").Value = 134
$ws.Range("J16").Value = 92
$ws.Range("K16").Value = 45.652173913043
$ws.Range("L16").Value = 41.052631578947
$ws.Range("M16").Value = 10.743801652892
$ws.Range("N16").Value = -81.046676096181

# ---- Row 17 (Fel. Assault) -------------------------------------------------
$ws.Range("C17").Value = 4
$ws.Range("E17").Value = -33.333333333333
$ws.Range("F17").Value = 22
$ws.Range("G17").Value = 18
$ws.Range("H17").Value = 22.222222222222
$ws.Range("I17").Value = 347
$ws.Range("J17").Value = 266
$ws.Range("K17").Value = 30.451127819548
$ws.Range("L17").Value = 47.659574468085
$ws.Range("M17").Value = 118.238993710692
$ws.Range("N17").Value = -33.652007648183

# ---- Row 18 (Burglary) -----------------------------------------------------
Set-Placeholder $ws.Range("C18") $zeroDonor
Set-Placeholder $ws.Range("D18") $zeroDonor
Set-Placeholder $ws.Range("E18") $naDonor
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 3
$ws.Range("H18").Value = 200
$ws.Range("L18").Value = 6.756756756756
$ws.Range("M18").Value = -41.911764705882
$ws.Range("N18").Value = -89.673202614379

# ---- Row 19 (Gr. Larceny) --------------------------------------------------
$ws.Range("C19").Value = 4
$ws.Range("D19").Value = 5
$ws.Range("E19").Value = -20
$ws.Range("F19").Value = 16
$ws.Range("G19").Value = 9
$ws.Range("H19").Value = 77.777777777777
$ws.Range("I19").Value = 180
$ws.Range("J19").Value = 180
$ws.Range("L19").Value = 10.429447852760
$ws.Range("M19").Value = 78.217821782178
$ws.Range("N19").Value = -44.954128440367

# ---- Row 20 (G.L.A.) --------------------------------------------------------
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = -83.333333333333
$ws.Range("G20").Value = 21
$ws.Range("H20").Value = -80.952380952380
$ws.Range("I20").Value = 75
$ws.Range("J20").Value = 82
$ws.Range("K20").Value = -8.536585365853
$ws.Range("L20").Value = 36.363636363636
$ws.Range("M20").Value = -3.846153846153
$ws.Range("N20").Value = -83.695652173913

# ---- Row 21 (TOTAL) ---------------------------------------------------------
$ws.Range("C21").Value = 11
$ws.Range("D21").Value = 19
$ws.Range("E21").Value = -42.105263157894
$ws.Range("F21").Value = 63
$ws.Range("G21").Value = 60
$ws.Range("H21").Value = 5
$ws.Range("I21").Value = 845
$ws.Range("J21").Value = 715
$ws.Range("K21").Value = 18.181818181818
$ws.Range("L21").Value = 30.602782071097
$ws.Range("M21").Value = 36.952998379254
$ws.Range("N21").Value = -70.246478873239

# ---- Row 22 (Transit) --------------------------------------------------------
Set-Placeholder $ws.Range("C22") $zeroDonor
$ws.Range("L22").Value = 10

# ---- Row 23 (Housing) --------------------------------------------------------
Set-Placeholder $ws.Range("C23") $zeroDonor
$ws.Range("E23").Value = -100
$ws.Range("F23").Value = 3
$ws.Range("H23").Value = -75
$ws.Range("J23").Value = 90
$ws.Range("K23").Value = -8.888888888888
$ws.Range("L23").Value = -5.747126436781

# ---- Row 24 (Petit Larceny) ---------------------------------------------------
$ws.Range("C24").Value = 5
$ws.Range("D24").Value = 13
$ws.Range("E24").Value = -61.538461538461
$ws.Range("F24").Value = 50
$ws.Range("G24").Value = 50
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 671
$ws.Range("J24").Value = 587
$ws.Range("K24").Value = 14.310051107325
$ws.Range("L24").Value = 14.700854700854
$ws.Range("M24").Value = 89.548022598870

# ---- Row 25 (Retail Theft) -----------------------------------------------------
Set-Placeholder $ws.Range("C25") $zeroDonor
Set-Placeholder $ws.Range("D25") $zeroDonor
Set-Placeholder $ws.Range("E25") $naDonor
$ws.Range("F25").Value = 1
$ws.Range("H25").Value = 0
$ws.Range("L25").Value = -38.461538461538

# ---- Row 26 (Misd. Assault) -----------------------------------------------------
$ws.Range("C26").Value = 9
$ws.Range("E26").Value = 12.5
$ws.Range("F26").Value = 42
$ws.Range("G26").Value = 26
$ws.Range("H26").Value = 61.538461538461
$ws.Range("I26").Value = 456
$ws.Range("J26").Value = 418
$ws.Range("K26").Value = 9.090909090909
$ws.Range("L26").Value = 22.252010723860
$ws.Range("M26").Value = 11.491442542787

# ---- Row 27 (UCR Rape*) -----------------------------------------------------------
Set-Placeholder $ws.Range("C27") $zeroDonor
Set-Placeholder $ws.Range("D27") $zeroDonor
Set-Placeholder $ws.Range("E27") $naDonor
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = -25

# ---- Row 28 (Other Sex Crimes) -----------------------------------------------------
# D28/E28/G28/H28 flip from the text "no data" placeholder style (13) to the normal
# numeric styles (14 / 15); copy format from a stable numeric donor first, then write the value.
$numDonor = $ws.Range("I30")   # style 14 (integer count format)
$pctDonor = $ws.Range("K30")   # style 15 (percentage format)

$numDonor.Copy($ws.Range("D28"))
$ws.Range("D28").Value = 1

$pctDonor.Copy($ws.Range("E28"))
$ws.Range("E28").Value = -100

$ws.Range("F28").Value = 4

$numDonor.Copy($ws.Range("G28"))
$ws.Range("G28").Value = 1

$pctDonor.Copy($ws.Range("H28"))
$ws.Range("H28").Value = 300

$ws.Range("J28").Value = 36
$ws.Range("K28").Value = 36.111111111111

# ---- Row 29 (Shooting Vic.) -----------------------------------------------------
$ws.Range("L29").Value = -33.333333333333

# ---- Row 30 (Shooting Inc.) -----------------------------------------------------
$ws.Range("L30").Value = -23.809523809523
